$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '69.140.52'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '3.741.57'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'601.60"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = "'167.39"
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').Value = '3.740.09'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').Value = "'0.170"
$ws.Range('E10').Value = '  +4.13%  '
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').Value = "'38.09"
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  +1.86%  '
$ws.Range('D15').Value = '4.368.49'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '3.746.31'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '69.059.19'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D21').Value = "'11.26"
$ws.Range('E21').Value = '  +11.75%  '
$ws.Range('D22').Value = "'492.62"
$ws.Range('E22').Value = '  -0.75%  '
$ws.Range('D23').Value = "'0.728"
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('E24').Value = '  +8.25%  '
$ws.Range('D25').Value = "'84.91"
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('E26').Value = '  -0.45%  '
$ws.Range('D27').Value = "'12.23"
$ws.Range('E27').Value = '  -0.93%  '
$ws.Range('D28').Value = "'10.06"
$ws.Range('E28').Value = '  -0.46%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = "'8.25"
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.97"
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').Value = "'31.52"
$ws.Range('E33').Value = '  -0.51%  '
$ws.Range('D34').Value = '3.887.61'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('D35').Value = '3.675.51'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E38').Value = '  +0.35%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "'0.139"
$ws.Range('E39').Value = '  +5.71%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = "'5.94"
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('D41').Value = "'0.326"
$ws.Range('E41').Value = '  +0.32%  '
$ws.Range('E42').Value = '  +6.24%  '
$ws.Range('D43').Value = "'48.81"
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('D44').Value = "'1.99"
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = "'423.20"
$ws.Range('E45').Value = '  -2.49%  '
$ws.Range('D46').Value = "'8.46"
$ws.Range('E46').Value = '  +0.61%  '
$ws.Range('E48').Value = '  -1.25%  '
$ws.Range('D49').Value = "'141.52"
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').Value = '2.781.15'
$ws.Range('E50').Value = '  +1.45%  '
$ws.Range('D51').Value = "'0.0353"
$ws.Range('E51').Value = '  +0.12%  '
